$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 2072.5454
$ws.Cells.Item(33, 9).Value = 1749.7142
$ws.Cells.Item(33, 11).Value = 1749.7142
$ws.Cells.Item(33, 13).Value = -1520.7142

$ws.Cells.Item(57, 8).Value = 53967
$ws.Cells.Item(57, 10).Value = 53967
$ws.Cells.Item(57, 12).Value = 161901
$ws.Cells.Item(57, 14).Value = -162899

$ws.Cells.Item(98, 8).Value = 983.6429000000001
$ws.Cells.Item(98, 9).Value = 363.44446
$ws.Cells.Item(98, 11).Value = 363.44446
$ws.Cells.Item(98, 13).Value = 1134.55554

$ws.Cells.Item(116, 8).Value = 9511.666999999999
$ws.Cells.Item(116, 9).Value = 11627.728
$ws.Cells.Item(116, 10).Value = 3692.5
$ws.Cells.Item(116, 11).Value = 11627.728
$ws.Cells.Item(116, 12).Value = 3692.5
$ws.Cells.Item(116, 13).Value = -8185.727999999999
$ws.Cells.Item(116, 14).Value = -10576.5

$ws.Cells.Item(122, 8).Value = 983.6429000000001
$ws.Cells.Item(122, 9).Value = 363.44446
$ws.Cells.Item(122, 11).Value = 1090.33338
$ws.Cells.Item(122, 13).Value = 1359.66662

$ws.Cells.Item(132, 8).Value = 6133.75
$ws.Cells.Item(132, 9).Value = 6453.5293
$ws.Cells.Item(132, 11).Value = 19360.5879
$ws.Cells.Item(132, 13).Value = -16830.5879

$ws.Cells.Item(138, 8).Value = 2336.4065
$ws.Cells.Item(138, 9).Value = 2449.6875
$ws.Cells.Item(138, 10).Value = 2312.24
$ws.Cells.Item(138, 11).Value = 7349.0625
$ws.Cells.Item(138, 12).Value = 6936.719999999999
$ws.Cells.Item(138, 13).Value = -2209.0625
$ws.Cells.Item(138, 14).Value = -17216.72

$ws.Cells.Item(141, 8).Value = 8158.2
$ws.Cells.Item(141, 9).Value = 3411.125
$ws.Cells.Item(141, 10).Value = 13583.429
$ws.Cells.Item(141, 11).Value = 10233.375
$ws.Cells.Item(141, 12).Value = 40750.287
$ws.Cells.Item(141, 13).Value = -5053.375
$ws.Cells.Item(141, 14).Value = -51110.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 430570.4
$ws.Cells.Item(32, 9).Value = 485248.34
$ws.Cells.Item(32, 10).Value = 23523.666
$ws.Cells.Item(32, 11).Value = 485248.34
$ws.Cells.Item(32, 12).Value = 23523.666
$ws.Cells.Item(32, 13).Value = -484961.34
$ws.Cells.Item(32, 14).Value = -24097.666

$ws.Cells.Item(52, 8).Value = 50000
$ws.Cells.Item(52, 10).Value = 50000
$ws.Cells.Item(52, 12).Value = 50000
$ws.Cells.Item(52, 14).Value = -50636

$ws.Cells.Item(61, 8).Value = 2498.52
$ws.Cells.Item(61, 9).Value = 1797.8235
$ws.Cells.Item(61, 10).Value = 3987.5
$ws.Cells.Item(61, 11).Value = 1797.8235
$ws.Cells.Item(61, 12).Value = 3987.5
$ws.Cells.Item(61, 13).Value = -1585.8235
$ws.Cells.Item(61, 14).Value = -4411.5

$ws.Cells.Item(92, 8).Value = 67498.336
$ws.Cells.Item(92, 10).Value = 67498.336
$ws.Cells.Item(92, 12).Value = 67498.336
$ws.Cells.Item(92, 14).Value = -72490.336

$ws.Cells.Item(136, 8).Value = 2498.52
$ws.Cells.Item(136, 9).Value = 1797.8235
$ws.Cells.Item(136, 10).Value = 3987.5
$ws.Cells.Item(136, 11).Value = 5393.470499999999
$ws.Cells.Item(136, 12).Value = 11962.5
$ws.Cells.Item(136, 13).Value = -2843.470499999999
$ws.Cells.Item(136, 14).Value = -17062.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 953.8421
$ws.Cells.Item(94, 9).Value = 744.5625
$ws.Cells.Item(94, 11).Value = 744.5625
$ws.Cells.Item(94, 13).Value = -293.5625

$ws.Cells.Item(99, 8).Value = 1316.0769
$ws.Cells.Item(99, 9).Value = 1180.9
$ws.Cells.Item(99, 11).Value = 1180.9
$ws.Cells.Item(99, 13).Value = 317.0999999999999

$ws.Cells.Item(100, 8).Value = 82596.664
$ws.Cells.Item(100, 10).Value = 82596.664
$ws.Cells.Item(100, 12).Value = 82596.664
$ws.Cells.Item(100, 14).Value = -84760.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 212.77777
$ws.Cells.Item(19, 9).Value = 212.77777
$ws.Cells.Item(19, 11).Value = 212.77777
$ws.Cells.Item(19, 13).Value = -42.77777

$ws.Cells.Item(24, 8).Value = 212.77777
$ws.Cells.Item(24, 9).Value = 212.77777
$ws.Cells.Item(24, 11).Value = 212.77777
$ws.Cells.Item(24, 13).Value = -42.77777

$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 8).Value = 33450
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 33450
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 14).Value = -35696

$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 8).Value = 33450
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 33450
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 14).Value = -111582

$ws.Cells.Item(87, 8).Value = 55165
$ws.Cells.Item(87, 10).Value = 55165
$ws.Cells.Item(87, 12).Value = 55165
$ws.Cells.Item(87, 14).Value = -57537

$ws.Cells.Item(90, 8).Value = 55165
$ws.Cells.Item(90, 10).Value = 55165
$ws.Cells.Item(90, 12).Value = 165495
$ws.Cells.Item(90, 14).Value = -177351

$ws.Cells.Item(99, 8).Value = 1867.9048
$ws.Cells.Item(99, 9).Value = 1804
$ws.Cells.Item(99, 10).Value = 1878.5555
$ws.Cells.Item(99, 11).Value = 1804
$ws.Cells.Item(99, 12).Value = 1878.5555
$ws.Cells.Item(99, 13).Value = -306
$ws.Cells.Item(99, 14).Value = -4874.5555

$ws.Cells.Item(107, 8).Value = 2718073
$ws.Cells.Item(107, 9).Value = 4464907.5
$ws.Cells.Item(107, 10).Value = 774.44446
$ws.Cells.Item(107, 11).Value = 4464907.5
$ws.Cells.Item(107, 12).Value = 774.44446
$ws.Cells.Item(107, 13).Value = -4462987.5
$ws.Cells.Item(107, 14).Value = -4614.44446

$ws.Cells.Item(126, 8).Value = 1867.9048
$ws.Cells.Item(126, 9).Value = 1804
$ws.Cells.Item(126, 10).Value = 1878.5555
$ws.Cells.Item(126, 11).Value = 5412
$ws.Cells.Item(126, 12).Value = 5635.666499999999
$ws.Cells.Item(126, 13).Value = -2942
$ws.Cells.Item(126, 14).Value = -10575.6665

$ws.Cells.Item(139, 8).Value = 30780
$ws.Cells.Item(139, 10).Value = 30780
$ws.Cells.Item(139, 12).Value = 30780
$ws.Cells.Item(139, 14).Value = -41060

$ws.Cells.Item(140, 8).Value = 37750
$ws.Cells.Item(140, 10).Value = 37750
$ws.Cells.Item(140, 12).Value = 37750
$ws.Cells.Item(140, 14).Value = -48110

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 559.41174
$ws.Cells.Item(5, 9).Value = 500.66666
$ws.Cells.Item(5, 11).Value = 1501.99998
$ws.Cells.Item(5, 13).Value = -1389.99998

$ws.Cells.Item(24, 8).Value = 2500
$ws.Cells.Item(24, 9).Value = 1000
$ws.Cells.Item(24, 10).Value = 4000
$ws.Cells.Item(24, 11).Value = 3000
$ws.Cells.Item(24, 12).Value = 12000
$ws.Cells.Item(24, 14).Value = -12460
$ws.Cells.Item(24, 13).Value = -2770

$ws.Cells.Item(135, 8).Value = 559.41174
$ws.Cells.Item(135, 9).Value = 500.66666
$ws.Cells.Item(135, 11).Value = 4505.99994
$ws.Cells.Item(135, 13).Value = -1970.99994

$ws.Cells.Item(137, 8).Value = 5215095.5
$ws.Cells.Item(137, 9).Value = 20848782
$ws.Cells.Item(137, 10).Value = 3867.0417
$ws.Cells.Item(137, 11).Value = 62546346
$ws.Cells.Item(137, 12).Value = 11601.1251
$ws.Cells.Item(137, 13).Value = -62541246
$ws.Cells.Item(137, 14).Value = -21801.1251

$ws.Cells.Item(140, 8).Value = 1625.75
$ws.Cells.Item(140, 10).Value = 2028.1111
$ws.Cells.Item(140, 12).Value = 6084.3333
$ws.Cells.Item(140, 14).Value = -16444.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 91.5
$ws.Cells.Item(2, 9).Value = 35
$ws.Cells.Item(2, 10).Value = 148
$ws.Cells.Item(2, 11).Value = 35
$ws.Cells.Item(2, 12).Value = 148
$ws.Cells.Item(2, 13).Value = 78
$ws.Cells.Item(2, 14).Value = -374

$ws.Cells.Item(15, 8).Value = 18833.334
$ws.Cells.Item(15, 10).Value = 18833.334
$ws.Cells.Item(15, 12).Value = 18833.334
$ws.Cells.Item(15, 14).Value = -19409.334

$ws.Cells.Item(81, 8).Value = 18833.334
$ws.Cells.Item(81, 10).Value = 18833.334
$ws.Cells.Item(81, 12).Value = 18833.334
$ws.Cells.Item(81, 14).Value = -20829.334

$ws.Cells.Item(84, 8).Value = 18833.334
$ws.Cells.Item(84, 10).Value = 18833.334
$ws.Cells.Item(84, 12).Value = 56500.00199999999
$ws.Cells.Item(84, 14).Value = -66484.00199999999

$ws.Cells.Item(97, 8).Value = 1178
$ws.Cells.Item(97, 9).Value = 722.5
$ws.Cells.Item(97, 11).Value = 722.5
$ws.Cells.Item(97, 13).Value = -226.5

$ws.Cells.Item(102, 8).Value = 912.0625
$ws.Cells.Item(102, 9).Value = 833.4545000000001
$ws.Cells.Item(102, 10).Value = 1085
$ws.Cells.Item(102, 11).Value = 833.4545000000001
$ws.Cells.Item(102, 12).Value = 1085
$ws.Cells.Item(102, 13).Value = 788.5454999999999
$ws.Cells.Item(102, 14).Value = -4329

$ws.Cells.Item(126, 8).Value = 2111.9443
$ws.Cells.Item(126, 9).Value = 1900.0769
$ws.Cells.Item(126, 11).Value = 5700.2307
$ws.Cells.Item(126, 13).Value = -3230.2307

$ws.Cells.Item(137, 8).Value = 20000
$ws.Cells.Item(137, 10).Value = 20000
$ws.Cells.Item(137, 12).Value = 20000
$ws.Cells.Item(137, 14).Value = -30200

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(43, 8).Value = 70000
$ws.Cells.Item(43, 10).Value = 70000
$ws.Cells.Item(43, 12).Value = 70000
$ws.Cells.Item(43, 14).Value = -70386

$ws.Cells.Item(132, 8).Value = 4362.278
$ws.Cells.Item(132, 10).Value = 4878
$ws.Cells.Item(132, 12).Value = 14634
$ws.Cells.Item(132, 14).Value = -19694

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(92, 8).Value = 50000
$ws.Cells.Item(92, 10).Value = 50000
$ws.Cells.Item(92, 12).Value = 50000
$ws.Cells.Item(92, 14).Value = -54992

$ws.Cells.Item(138, 8).Value = 99429
$ws.Cells.Item(138, 10).Value = 99429
$ws.Cells.Item(138, 12).Value = 99429
$ws.Cells.Item(138, 14).Value = -109709
